$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.532163619995117
$ws.Range("B1").Value = 4.634356021881104
$ws.Range("C1").Value = 8.386774063110352
$ws.Range("D1").Value = 8.699810981750488
$ws.Range("E1").Value = 5.106639862060547
